# The workbook tracks weekly Brócoli wholesale-market prices. A new week
# of data (2022-05-13, date serial 44694) is inserted at the top of the
# "historical" block (rows 640-641), pushing all subsequent weeks down by
# two rows. Everything else (Mercado/Región/Categoría/etc.) for those two
# new rows mirrors the row immediately below them (same market, same
# product, "Primera" / "Segunda" quality pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 640 - this shifts the existing rows
# 640:757 down to 642:759 and carries their formatting (incl. the date
# number format on column D) along with them.
$ws.Rows("640:641").Insert()

# New row 640: Brócoli, Primera
$ws.Range("A640").Value = 6
$ws.Range("B640").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C640").Value = "Metropolitana"
$ws.Range("D640").Value = 44694
$ws.Range("E640").Value = 13
$ws.Range("F640").Value = 100112023
$ws.Range("G640").Value = "Brócoli"
$ws.Range("H640").Value = "Sin especificar"
$ws.Range("I640").Value = "Primera"
$ws.Range("J640").Value = 15300
$ws.Range("K640").Value = 550
$ws.Range("L640").Value = 700
$ws.Range("M640").Value = 628
$ws.Range("N640").Value = "$/unidad"
$ws.Range("O640").Value = "Región Metropolitana"
$ws.Range("P640").Value = 628
$ws.Range("Q640").Value = 1
$ws.Range("R640").Value = "Hortaliza"

# New row 641: Brócoli, Segunda
$ws.Range("A641").Value = 6
$ws.Range("B641").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C641").Value = "Metropolitana"
$ws.Range("D641").Value = 44694
$ws.Range("E641").Value = 13
$ws.Range("F641").Value = 100112023
$ws.Range("G641").Value = "Brócoli"
$ws.Range("H641").Value = "Sin especificar"
$ws.Range("I641").Value = "Segunda"
$ws.Range("J641").Value = 6200
$ws.Range("K641").Value = 400
$ws.Range("L641").Value = 500
$ws.Range("M641").Value = 453
$ws.Range("N641").Value = "$/unidad"
$ws.Range("O641").Value = "Región Metropolitana"
$ws.Range("P641").Value = 453
$ws.Range("Q641").Value = 1
$ws.Range("R641").Value = "Hortaliza"
